$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
# Columns A and E hold numeric-looking text ("-204109293", "171219") that must be
# stored as text (shared string), not auto-converted to a number, so force the
# cell's number format to Text ("@") right before assigning the value.
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "-204109293"

$ws.Range("B4").Value = "asd"
$ws.Range("C4").Value = "asd"
$ws.Range("D4").Value = "asd"

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "171219"

# Row 5
$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "-1715655579"

$ws.Range("B5").Value = "asd"
$ws.Range("C5").Value = "asdasd"
$ws.Range("D5").Value = "asd"

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "171219"
